# Regenerate merged AHB files
# Rename header columns from _old/_new to _FV2210/_FV2304, add a table
# over the data range, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1) -----------------------------------
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $text = $cell.Value()
    if ($text -like "*_old") {
        $cell.Value = ($text -replace "_old$", "_FV2210")
    } elseif ($text -like "*_new") {
        $cell.Value = ($text -replace "_new$", "_FV2304")
    }
}

# --- 2. Create a table (ListObject) over the used data range ----------
$rng = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add(1, $rng, [Type]::Missing, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
